$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell -> new text value. Values are quote-prefixed so Excel stores them
# as literal text (matching the source workbook, which keeps these as
# inline strings) instead of auto-converting numeric-looking text to numbers.
$changes = [ordered]@{
    'D2' = "'307.59"
    'E2' = "'-2.16%"
    'D3' = "'37.61"
    'E3' = "'-4.81%"
    'D4' = "'5.101"
    'E4' = "'-0.28%"
    'D5' = "'0.07880"
    'E5' = "'-3.76%"
    'D6' = "'1.980"
    'E6' = "'-3.48%"
    'D7' = "'4.334"
    'E7' = "'1.36%"
    'D8' = "'8.235"
    'E8' = "'-0.14%"
    'E9' = "'-6.51%"
    'D10' = "'0.9308"
    'E10' = "'-0.22%"
    'D11' = "'0.1302"
    'E11' = "'-7.70%"
    'D12' = "'0.1899"
    'E12' = "'-4.55%"
    'D13' = "'0.08879"
    'E13' = "'-2.87%"
    'D14' = "'0.03437"
    'E14' = "'-2.63%"
    'D15' = "'0.09744"
    'E15' = "'-0.72%"
    'D16' = "'0.001391"
    'E16' = "'-0.37%"
    'D17' = "'0.005877"
    'E17' = "'-6.25%"
    'E18' = "'1,774.23%"
    'D19' = "'3.563"
    'E19' = "'-2.67%"
    'D20' = "'0.3431"
    'E20' = "'-0.81%"
    'D21' = "'0.1288"
    'E21' = "'-1.26%"
    'E22' = "'1.93%"
    'E23' = "'1.40%"
    'D24' = "'0.04307"
    'E24' = "'-0.53%"
    'E25' = "'-0.58%"
    'D26' = "'0.004609"
    'E26' = "'-3.54%"
    'E27' = "'176.03%"
    'D39' = "'0.02342"
    'E39' = "'4.86%"
    'D40' = "'0.05036"
    'E40' = "'-3.98%"
    'D41' = "'0.007533"
    'E41' = "'0.24%"
    'D42' = "'0.009766"
    'E42' = "'-0.23%"
    'D43' = "'0.1355"
    'E43' = "'-1.65%"
    'D44' = "'0.002089"
    'E44' = "'-1.38%"
    'D45' = "'0.008010"
    'E45' = "'-15.44%"
    'D46' = "'0.00006523"
    'E46' = "'0.95%"
    'E47' = "'-0.16%"
    'D48' = "'0.002997"
    'E48' = "'8.22%"
    'E50' = "'-0.16%"
    'E51' = "'-0.16%"
}

foreach ($cellRef in $changes.Keys) {
    $ws.Range($cellRef).Value = $changes[$cellRef]
}
